$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Personnes" (sheet1): fix a couple of bios, and add new contributors
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Personnes")

# Row 4 (Marcela Gori): trim the bio text
$ws.Cells.Item(4, 6).Value = "Élue locale à Anderlecht, néo-libérale"

# Row 5 (Daniel Rodenstein): trim the bio text (keep trailing newline)
$ws.Cells.Item(5, 6).Value = "Membre de l'Institut Jonathas`n"

# New row 6: Michel Cristt, Bruxelles
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = "Michel"
$ws.Cells.Item(6, 3).Value = "Cristt"
$ws.Cells.Item(6, 4).Value = "Bruxelles"
$ws.Cells.Item(6, 6).Value = "Momomo`n"

# New row 7: Cat Firch, Bruxelles
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = "Cat"
$ws.Cells.Item(7, 3).Value = "Firch"
$ws.Cells.Item(7, 4).Value = "Bruxelles"
$ws.Cells.Item(7, 6).Value = "Nininiin"

# New rows 8-10: placeholder IDs only
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(10, 1).Value = 9

# Former row 6 (Lucas Ablotia) shifts down to row 11
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "Lucas Ablotia"
$ws.Cells.Item(11, 6).Value = "Sionniste, lien avec Belgian Friends of Israel, islamophobe, pro-Trump"
$ws.Cells.Item(6, 2).Value = "Michel"

# Clear the old row-6 leftovers that are no longer part of that row
$ws.Cells.Item(6, 6).Value = "Momomo`n"

# The multi-line bio cells get a (no-op) wrap-text toggle, matching the new
# style entry introduced in the workbook
$ws.Range("F5").WrapText = $true
$ws.Range("F5").WrapText = $false
$ws.Range("F6").Style = $ws.Range("F5").Style

$ws.Activate()
$ws.Range("G12").Select()

# ---------------------------------------------------------------------------
# Sheet "Articles" (sheet2): add article rows, one with a hyperlink
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Articles")

$ws2.Cells.Item(2, 1).Value = 1
$ws2.Cells.Item(2, 3).Value = 1

$ws2.Cells.Item(3, 1).Value = 2
$ws2.Cells.Item(3, 2).Value = "Progressisme et deconstructrion"
$ws2.Cells.Item(3, 3).Value = 1
$ws2.Hyperlinks.Add($ws2.Range("E3"), "https://www.21news.be/progressisme-et-deconstruction-les-racines-ideologiques-du-mal-etre-contemporain/")
$ws2.Range("E3").Style = "Lien hypertexte"

$ws2.Cells.Item(4, 1).Value = 3
$ws2.Cells.Item(4, 3).Value = 2

$ws2.Cells.Item(5, 1).Value = 4
$ws2.Cells.Item(5, 3).Value = 3

$ws2.Activate()
$ws2.Range("C4").Select()

# ---------------------------------------------------------------------------
# Sheet "Relations" (sheet3): no content changes, just move the selection
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Relations")
$ws3.Activate()
$ws3.Range("E36").Select()

# ---------------------------------------------------------------------------
# Sheet "Debunks" (sheet4): populate with header + data, one hyperlink
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Debunks")

$ws4.Cells.Item(1, 1).Value = "Article_ID"
$ws4.Cells.Item(1, 2).Value = "Article_Debunk_ID"
$ws4.Cells.Item(1, 3).Value = "Type"

$ws4.Cells.Item(2, 1).Value = 2
$ws4.Cells.Item(2, 2).Value = 1
$ws4.Hyperlinks.Add($ws4.Range("C2"), "https://www.arretsurimages.net/articles/marie-estelle-dupont-psycho-couacs-a-droite-toute")
$ws4.Range("C2").Style = "Lien hypertexte"

$ws4.Cells.Item(3, 2).Value = 2
$ws4.Cells.Item(4, 2).Value = 3
$ws4.Cells.Item(5, 2).Value = 4
$ws4.Cells.Item(6, 2).Value = 5

$ws4.Activate()
$ws4.Range("C16").Select()

# ---------------------------------------------------------------------------
# Leave "Personnes" as the active sheet / tab, matching the saved selection
# ---------------------------------------------------------------------------
$ws.Activate()
